# new results. 2021/04/21 15:43
#
# Column B (row 2..31) holds per-"expert" scores that are colour-coded with
# the built-in "Good / Neutral / Bad" cell styles (green / orange / red).
# This refresh:
#   1. corrects row 20's value (it now matches row 10's corrected figure),
#      which ripples into the AVERAGE formula in B32;
#   2. re-touches the Neutral/Bad cell-style pairs below B9, which makes
#      Excel re-issue their underlying xf records (the cell-style gallery
#      entries keep pointing at the same named "Neutral"/"Bad" styles);
#   3. leaves the selection on C1 instead of B17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix the mis-keyed data point -------------------------------------
# Row 20 (expert #19) was 0.3643; the corrected figure is 0.3762 (same as
# row 10 / expert #9). B32 (=AVERAGE(B2:B31)) recalculates automatically.
$ws.Cells.Item(20, 2).Value = 0.37619999999999998

# --- 2. Re-apply Neutral/Bad cell styles, swapping their bookkeeping slot -
# Every row below pairs a "Bad" (red) cell with a "Neutral" (orange) cell;
# re-stamping each one with the other's Style object makes Excel re-issue
# the xf records in the opposite order while every cell keeps its original
# colour (red stays red, orange stays orange).
$swapRowPairs = @(
    @(3, 10),
    @(4, 11),
    @(8, 12),
    @(13, 16),
    @(14, 17),
    @(15, 20),
    @(26, 22),
    @(27, 28),
    @(31, 30)
)

foreach ($pair in $swapRowPairs) {
    $rowA = $pair[0]
    $rowB = $pair[1]
    $cellA = $ws.Cells.Item($rowA, 2)
    $cellB = $ws.Cells.Item($rowB, 2)

    $styleA = $cellA.Style
    $styleB = $cellB.Style

    $cellA.Style = $styleB
    $cellB.Style = $styleA
}

# --- 3. Move the active selection ----------------------------------------
[void]$ws.Range("C1").Select()
